$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23 ---
$ws.Range("A23").Value = 20
$ws.Range("B23").Value = "Generate UI (analytics, solutions, solcategoties, tickets)"
$ws.Range("D23").Value = "Task 17 is finished"

# --- Fill Done-criteria column out of row order (matches original authoring order) ---
$ws.Range("D25").Value = "Task 16 is finished"
$ws.Range("D24").Value = "Task 13 is finished"
$ws.Range("D26").Value = "Task 14 is finished"
$ws.Range("D27").Value = "Task 15 is finished"

# --- Fill Description column for the remaining new rows ---
$ws.Range("B24").Value = "Start developing Analytics app UI part"
$ws.Range("B25").Value = "Start developing Solution Categories app UI part"
$ws.Range("B26").Value = "Start developing solutions app UI part"
$ws.Range("B27").Value = "Start developing tickets app UI part"

# --- Remaining cells reusing already-existing shared strings ---
$ws.Range("C23").Value = "Shamil"
$ws.Range("E23").Value = "test"
$ws.Range("F23").Value = "Done"
$ws.Range("G23").Value = "1 hour"

$ws.Range("A24").Value = 21
$ws.Range("C24").Value = "Shamil"
$ws.Range("F24").Value = "New"

$ws.Range("A25").Value = 22
$ws.Range("C25").Value = "Marcia"
$ws.Range("F25").Value = "New"

$ws.Range("A26").Value = 23
$ws.Range("C26").Value = "Ilkay"
$ws.Range("F26").Value = "New"

$ws.Range("A27").Value = 24
$ws.Range("C27").Value = "Marcus"
$ws.Range("F27").Value = "New"

# --- Formatting: columns A:E for rows 23-27 get wrap-text style (same as rest of table) ---
$ws.Range("A23:E23").WrapText = $true
$ws.Range("A24:D24").WrapText = $true
$ws.Range("A25:D25").WrapText = $true
$ws.Range("A26:D26").WrapText = $true
$ws.Range("A27:D27").WrapText = $true

# --- F23 (status "Done") gets the green fill + wrap-text style used elsewhere for Done ---
$ws.Range("F23").Interior.Color = 5296274
$ws.Range("F23").WrapText = $true

# --- Selection matches the last-edited cell ---
$ws.Range("B24").Select()
